$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 6762
$ws1.Range("F7").Value = 7
$ws1.Range("F10").Value = 6330
$ws1.Range("F15").Value = 104
$ws1.Range("F19").Value = 373
$ws1.Range("F22").Value = 4695
$ws1.Range("F23").Value = 69
$ws1.Range("F25").Value = 147
$ws1.Range("F27").Value = 94

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 6762
$ws4.Range("F7").Value = 7
$ws4.Range("F10").Value = 6330
$ws4.Range("F15").Value = 104
$ws4.Range("F19").Value = 373
$ws4.Range("F22").Value = 4695
$ws4.Range("F24").Value = 69
$ws4.Range("F26").Value = 147
$ws4.Range("F28").Value = 94
